$d = $word.ActiveDocument
$d.Content.Find.Execute("MacDonald’s Seating Area", $true, $false, $false, $false, $false,
                         $true, 1, $false, "MacDonald’s Seating", 2)
